$wb = $excel.ActiveWorkbook

# Sheet "展览": F4 850 -> 853, F5 74 -> 75
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 853
$wsExhibit.Range("F5").Value = 75

# Sheet "全部类型": F5 850 -> 853, F6 74 -> 75
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 853
$wsAll.Range("F6").Value = 75

$wb.Save()
